# Update the build timestamp embedded in various text cells to reflect the
# new release build time (January 30 2026 16.19.47 EST -> February 02 2026 12.49.33 EST).

$wb = $excel.ActiveWorkbook

$oldStamp = "January 30 2026 16.19.47 EST"
$newStamp = "February 02 2026 12.49.33 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet: version banner (A2) and citation text (A6) ---
$a2Text = $wsAbout.Range("A2").Value()
$a2Text = $a2Text.Replace($oldStamp, $newStamp)
$wsAbout.Range("A2").Value = $a2Text

$a6Text = $wsAbout.Range("A6").Value()
$a6Text = $a6Text.Replace($oldStamp, $newStamp)
$wsAbout.Range("A6").Value = $a6Text

# --- "Boundaries and methane sources" sheet: build_version column (S2:S8) ---
for ($row = 2; $row -le 8; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = 19
    $cellText = $cell.Value()
    $cellText = $cellText.Replace($oldStamp, $newStamp)
    $cell.Value = $cellText
}
